$d = $word.ActiveDocument
$vtab = [char]11

function Set-ParaText($paraIndex, $newText) {
    $p = $d.Paragraphs($paraIndex).Range
    $r = $d.Range($p.Start, $p.End - 1)
    $r.Text = $newText
}

# Paragraph 1: Title
Set-ParaText 1 "Exploring the Realm of Science: Unveiling the Secrets of Nature"

# Paragraph 2: Author name
Set-ParaText 2 "Dr. Sofia Richards"

# Paragraph 3: Email
Set-ParaText 3 "sofiarichards@scienceacademy.edu"

# Paragraph 5: Main body paragraph (four "block" sections separated by double line-breaks)
$blockA = "Science, an awe-inspiring realm of discovery, invites us to embark on an extraordinary journey into the profound mysteries of the universe. Through the rigorous study of science, we develop critical thinking skills, nurture analytical minds, and unlock the secrets of the natural world. Its meticulous methods allow us to unravel the enigmas that shape our existence"
$blockB = "Step into the fascinating world of chemistry, where atoms and molecules dance in a delicate symphony of reactions, revealing the intricate tapestry of matter. Witness the transformative power of chemical equations as substances undergo mesmerizing transformations, fostering an understanding of the transformative forces that govern the world around us"
$blockC = "Venture into the realm of biology, a symphony of life brimming with diversity and complexity. Explore the intricate mechanisms of cellular processes, unravel the genetic code that holds the blueprint of life, and marvel at the delicate balance of ecosystems. Embark on a microscopic odyssey, delving into the depths of DNA, the enigmatic blueprint that orchestrates the symphony of life"
$blockD = "Journey into the realm of medicine, a noble pursuit dedicated to alleviating human suffering. Discover the intricate workings of the human body, witness the body's remarkable ability to heal, and explore the frontiers of medical advancements. Learn about the selfless contributions of healthcare professionals, whose tireless efforts bring solace to those in need."

$para5 = $blockA + $vtab + $vtab + $blockB + $vtab + $vtab + $blockC + $vtab + $vtab + $blockD
Set-ParaText 5 $para5

# Paragraph 7: Summary body paragraph
# A trailing carriage-return character is appended so a brand-new, completely
# empty paragraph (no run at all, matching a bare <w:p/>) is created right
# after it -- mirroring the added <w:p/> at the end of the document body.
$cr = [char]13
$summary = "Science, in its myriad forms, offers a portal into the wonders of the universe. Through chemistry, we unveil the secrets of matter, unraveling the intricate dance of atoms and molecules. Biology unveils the symphony of life, revealing the delicate balance of ecosystems and the intricate workings of cellular processes. Medicine, a noble pursuit, alleviates human suffering and pushes the boundaries of medical knowledge. Science empowers us to comprehend the cosmos, unlock the mysteries of nature, and harness knowledge for the betterment of humanity." + $cr
Set-ParaText 7 $summary

Write-Output "Edits applied."
